# Restore C10 on the "Rules" sheet back to its earlier value of 1
# (was showing 18 before this revision-restore save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
